# Export with no is_pref and no lev distance
#
# The previous export computed `id` (column B) for "non-preferred"
# speaker-variant rows by finding the nearest *preferred* variant
# (Levenshtein distance) and reusing its slug, and it flagged the
# preferred rows with an "x" in the `is_prefered` column (D).
#
# The new export drops both behaviours:
#   - `id` (column B) is always just "#" + lowercase(speaker_variant)
#   - `is_prefered` (column D) is no longer populated, and the rows are
#     written out in a different (new) order reflecting the new logic.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New speaker_variant (column C) values, in final row order (rows 2-30).
$values = @(
    "Dah",
    "Queis",
    "Celia",
    "Thi",
    "Capj",
    "Fl",
    "Sil",
    "Osm,",
    "Celi",
    "Phafa",
    "Capi",
    "Thioss",
    "Phaf",
    "Janit",
    "Queissi",
    "Silactar",
    "Sultana Thiossem",
    "Osman",
    "Silact",
    "Osm",
    "Flora",
    "Dahout",
    "Ianit",
    "Thios",
    "Must",
    "Boust",
    "Ian",
    "Flo",
    "Cel"
)

$row = 2
foreach ($speakerVariant in $values) {
    $id = "#" + $speakerVariant.ToLower().Replace(" ", "-")

    $ws.Cells.Item($row, 2).Value = $id              # column B - id
    $ws.Cells.Item($row, 3).Value = $speakerVariant   # column C - speaker_variant
    $ws.Cells.Item($row, 4).Value = ""                # column D - is_prefered (cleared)

    $row = $row + 1
}
